$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data (columns D, L, M, N, O, P, Q, R, S, T) between row 3 <-> row 5
# and row 4 <-> row 6, leaving columns A, B, C, E, F, G, H, I, J, K untouched.

$ws.Range("D3").Value = 44334
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 12500
$ws.Range("Q3").Value = "$/caja 12 kilos empedrada"
$ws.Range("S3").Value = 1042
$ws.Range("T3").Value = 12

$ws.Range("D4").Value = 44330
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 15500
$ws.Range("Q4").Value = "$/caja 18 kilos granel"
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 861
$ws.Range("T4").Value = 18

$ws.Range("D5").Value = 44316
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 17500
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17750
$ws.Range("Q5").Value = "$/caja 16 kilos granel"
$ws.Range("S5").Value = 1109
$ws.Range("T5").Value = 16

$ws.Range("D6").Value = 44316
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("Q6").Value = "$/caja 16 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 16
